$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "NBA, Monday 12th Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

# Matchup / prediction / actual data rows (2-11)
$data = @(
    @("Indiana Pacers (30-25) vs Charlotte Hornets (11-41)", "Indiana Pacers (60.71%)", "Indiana Pacers (82.3%)"),
    @("Philadelphia 76ers (31-21) vs Cleveland Cavaliers (35-16)", "Cleveland Cavaliers (69.23%)", "Cleveland Cavaliers (86.1%)"),
    @("Chicago Bulls (25-28) vs Atlanta Hawks (24-29)", "Chicago Bulls (55.56%)", "Atlanta Hawks (62.6%)"),
    @("San Antonio Spurs (10-43) vs Toronto Raptors (19-34)", "Toronto Raptors (44.00%)", "Toronto Raptors (76.3%)"),
    @("New York Knicks (33-20) vs Houston Rockets (23-29)", "New York Knicks (70.37%)", "New York Knicks (54.9%)"),
    @("New Orleans Pelicans (31-22) vs Memphis Grizzlies (18-35)", "New Orleans Pelicans (60.00%)", "New Orleans Pelicans (77.0%)"),
    @("Denver Nuggets (36-17) vs Milwaukee Bucks (34-19)", "Denver Nuggets (84.00%)", "Milwaukee Bucks (53.6%)"),
    @("Washington Wizards (9-43) vs Dallas Mavericks (30-23)", "Dallas Mavericks (53.57%)", "Dallas Mavericks (83.8%)"),
    @("Golden State Warriors (25-25) vs Utah Jazz (26-27)", "Utah Jazz (70.83%)", "Utah Jazz (53.5%)"),
    @("Minnesota Timberwolves (36-16) vs Los Angeles Clippers (35-16)", "Los Angeles Clippers (80.00%)", "Los Angeles Clippers (64.7%)")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Move selection to the last populated cell, like the saved file shows
$ws.Range("A11").Select() | Out-Null
